$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAR")

$ws.Cells.Item(20, 1).Value = "Cakmakci"
$ws.Cells.Item(20, 2).Value = 2021
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = "Medium-textured"
$ws.Cells.Item(20, 5).Value = 2.450268029
$ws.Cells.Item(20, 6).Value = 0.0216
$ws.Cells.Item(20, 7).Value = 3
$ws.Cells.Item(20, 8).Value = 2.273155622
$ws.Cells.Item(20, 9).Value = 0.01296
$ws.Cells.Item(20, 10).Value = 3

$ws.Cells.Item(21, 1).Value = "Cakmakci"
$ws.Cells.Item(21, 2).Value = 2021
$ws.Cells.Item(21, 3).Value = 2
$ws.Cells.Item(21, 4).Value = "Medium-textured"
$ws.Cells.Item(21, 5).Value = 2.237841577
$ws.Cells.Item(21, 6).Value = 0.01584
$ws.Cells.Item(21, 7).Value = 3
$ws.Cells.Item(21, 8).Value = 2.125078392
$ws.Cells.Item(21, 9).Value = 0.03024
$ws.Cells.Item(21, 10).Value = 3

$ws.Cells.Item(22, 1).Value = "Cakmakci"
$ws.Cells.Item(22, 2).Value = 2021
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 4).Value = "Medium-textured"
$ws.Cells.Item(22, 5).Value = 1.83779367
$ws.Cells.Item(22, 6).Value = 0.04824
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(22, 8).Value = 1.782016545
$ws.Cells.Item(22, 9).Value = 0.02088
$ws.Cells.Item(22, 10).Value = 3

$ws.Cells.Item(23, 1).Value = "Cakmakci"
$ws.Cells.Item(23, 2).Value = 2021
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = "Medium-textured"
$ws.Cells.Item(23, 5).Value = 2.542725171
$ws.Cells.Item(23, 6).Value = 0.00864
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 2.592611908
$ws.Cells.Item(23, 9).Value = 0.04464
$ws.Cells.Item(23, 10).Value = 3

$ws.Cells.Item(24, 1).Value = "Cakmakci"
$ws.Cells.Item(24, 2).Value = 2021
$ws.Cells.Item(24, 3).Value = 2
$ws.Cells.Item(24, 4).Value = "Medium-textured"
$ws.Cells.Item(24, 5).Value = 2.379314409
$ws.Cells.Item(24, 6).Value = 0.0684
$ws.Cells.Item(24, 7).Value = 3
$ws.Cells.Item(24, 8).Value = 2.273155622
$ws.Cells.Item(24, 9).Value = 0.03744
$ws.Cells.Item(24, 10).Value = 3

$ws.Cells.Item(25, 1).Value = "Cakmakci"
$ws.Cells.Item(25, 2).Value = 2021
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 4).Value = "Medium-textured"
$ws.Cells.Item(25, 5).Value = 2.209616276
$ws.Cells.Item(25, 6).Value = 0.0108
$ws.Cells.Item(25, 7).Value = 3
$ws.Cells.Item(25, 8).Value = 2.019696257
$ws.Cells.Item(25, 9).Value = 0.04752
$ws.Cells.Item(25, 10).Value = 3

$ws.Cells.Item(26, 1).Value = "Cakmakci"
$ws.Cells.Item(26, 2).Value = 2021
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 4).Value = "Medium-textured"
$ws.Cells.Item(26, 5).Value = 2.287291337
$ws.Cells.Item(26, 6).Value = 0.03024
$ws.Cells.Item(26, 7).Value = 3
$ws.Cells.Item(26, 8).Value = 2.29436136
$ws.Cells.Item(26, 9).Value = 0.03096
$ws.Cells.Item(26, 10).Value = 3

$ws.Cells.Item(27, 1).Value = "Cakmakci"
$ws.Cells.Item(27, 2).Value = 2021
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = "Medium-textured"
$ws.Cells.Item(27, 5).Value = 2.223726047
$ws.Cells.Item(27, 6).Value = 0.01872
$ws.Cells.Item(27, 7).Value = 3
$ws.Cells.Item(27, 8).Value = 2.216670442
$ws.Cells.Item(27, 9).Value = 0.02736
$ws.Cells.Item(27, 10).Value = 3

$ws.Cells.Item(28, 1).Value = "Cakmakci"
$ws.Cells.Item(28, 2).Value = 2021
$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(28, 4).Value = "Medium-textured"
$ws.Cells.Item(28, 5).Value = 1.782016545
$ws.Cells.Item(28, 6).Value = 0.036
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(28, 8).Value = 2.103976207
$ws.Cells.Item(28, 9).Value = 0.018
$ws.Cells.Item(28, 10).Value = 3

$ws.Cells.Item(29, 1).Value = "Demir"
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = 2
$ws.Cells.Item(29, 4).Value = "Fine-textured"
$ws.Cells.Item(29, 5).Value = 2.379314409
$ws.Cells.Item(29, 6).Value = 0.072
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(29, 8).Value = 2.012682214
$ws.Cells.Item(29, 9).Value = 0.0648
$ws.Cells.Item(29, 10).Value = 3

$ws.Cells.Item(30, 1).Value = "Demir"
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = 2
$ws.Cells.Item(30, 4).Value = "Fine-textured"
$ws.Cells.Item(30, 5).Value = 2.259025679
$ws.Cells.Item(30, 6).Value = 0.108
$ws.Cells.Item(30, 7).Value = 3
$ws.Cells.Item(30, 8).Value = 1.956621193
$ws.Cells.Item(30, 9).Value = 0.0576
$ws.Cells.Item(30, 10).Value = 3

$ws.Cells.Item(31, 1).Value = "Demir"
$ws.Cells.Item(31, 2).Value = 2019
$ws.Cells.Item(31, 3).Value = 2
$ws.Cells.Item(31, 4).Value = "Fine-textured"
$ws.Cells.Item(31, 5).Value = 2.103976207
$ws.Cells.Item(31, 6).Value = 0.0792
$ws.Cells.Item(31, 7).Value = 3
$ws.Cells.Item(31, 8).Value = 2.047766703
$ws.Cells.Item(31, 9).Value = 0.1656
$ws.Cells.Item(31, 10).Value = 3

$ws.Cells.Item(32, 1).Value = "Demir"
$ws.Cells.Item(32, 2).Value = 2019
$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 4).Value = "Fine-textured"
$ws.Cells.Item(32, 5).Value = 2.19551226
$ws.Cells.Item(32, 6).Value = 0.1224
$ws.Cells.Item(32, 7).Value = 3
$ws.Cells.Item(32, 8).Value = 2.139153682
$ws.Cells.Item(32, 9).Value = 0.864
$ws.Cells.Item(32, 10).Value = 3

$ws.Cells.Item(33, 1).Value = "Demir"
$ws.Cells.Item(33, 2).Value = 2019
$ws.Cells.Item(33, 3).Value = 2
$ws.Cells.Item(33, 4).Value = "Fine-textured"
$ws.Cells.Item(33, 5).Value = 2.068834539
$ws.Cells.Item(33, 6).Value = 0.1224
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(33, 8).Value = 2.07586001
$ws.Cells.Item(33, 9).Value = 0.1008
$ws.Cells.Item(33, 10).Value = 3

$ws.Cells.Item(34, 1).Value = "Tunc"
$ws.Cells.Item(34, 2).Value = 2015
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = "Medium-textured"
$ws.Cells.Item(34, 5).Value = 2.280222758
$ws.Cells.Item(34, 6).Value = 0.162119956
$ws.Cells.Item(34, 7).Value = 3
$ws.Cells.Item(34, 8).Value = 1.823840902
$ws.Cells.Item(34, 9).Value = 0.399064506
$ws.Cells.Item(34, 10).Value = 3

$ws.Cells.Item(35, 1).Value = "Tunc"
$ws.Cells.Item(35, 2).Value = 2015
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 4).Value = "Medium-textured"
$ws.Cells.Item(35, 5).Value = 2.400585255
$ws.Cells.Item(35, 6).Value = 0.249415317
$ws.Cells.Item(35, 7).Value = 3
$ws.Cells.Item(35, 8).Value = 2.528484962
$ws.Cells.Item(35, 9).Value = 0.399064506
$ws.Cells.Item(35, 10).Value = 3

$ws.Cells.Item(36, 1).Value = "Tunc"
$ws.Cells.Item(36, 2).Value = 2015
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = "Medium-textured"
$ws.Cells.Item(36, 5).Value = 1.97763339
$ws.Cells.Item(36, 6).Value = 0.087295361
$ws.Cells.Item(36, 7).Value = 3
$ws.Cells.Item(36, 8).Value = 3.247719436
$ws.Cells.Item(36, 9).Value = 0.399064506
$ws.Cells.Item(36, 10).Value = 3

$ws.Cells.Item(37, 1).Value = "Tunc"
$ws.Cells.Item(37, 2).Value = 2015
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(37, 4).Value = "Medium-textured"
$ws.Cells.Item(37, 5).Value = 2.29436136
$ws.Cells.Item(37, 6).Value = 0.212003019
$ws.Cells.Item(37, 7).Value = 3
$ws.Cells.Item(37, 8).Value = 1.545969547
$ws.Cells.Item(37, 9).Value = 0.14964919
$ws.Cells.Item(37, 10).Value = 3

$ws.Cells.Item(38, 1).Value = "Tunc"
$ws.Cells.Item(38, 2).Value = 2015
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = "Medium-textured"
$ws.Cells.Item(38, 5).Value = 2.230783092
$ws.Cells.Item(38, 6).Value = 0.212003019
$ws.Cells.Item(38, 7).Value = 3
$ws.Cells.Item(38, 8).Value = 1.545969547
$ws.Cells.Item(38, 9).Value = 0.14964919
$ws.Cells.Item(38, 10).Value = 3

$ws.Cells.Item(39, 1).Value = "Tunc"
$ws.Cells.Item(39, 2).Value = 2015
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 4).Value = "Medium-textured"
$ws.Cells.Item(39, 5).Value = 1.782016545
$ws.Cells.Item(39, 6).Value = 0.14964919
$ws.Cells.Item(39, 7).Value = 3
$ws.Cells.Item(39, 8).Value = 1.545969547
$ws.Cells.Item(39, 9).Value = 0.14964919
$ws.Cells.Item(39, 10).Value = 3

$ws.Cells.Item(40, 1).Value = "Levy"
$ws.Cells.Item(40, 2).Value = 2014
$ws.Cells.Item(40, 3).Value = 2
$ws.Cells.Item(40, 4).Value = "Coarse-textured"
$ws.Cells.Item(40, 5).Value = 2.807226058
$ws.Cells.Item(40, 6).Value = 0.36
$ws.Cells.Item(40, 7).Value = 3
$ws.Cells.Item(40, 8).Value = 2.450268029
$ws.Cells.Item(40, 9).Value = 0.144
$ws.Cells.Item(40, 10).Value = 3

$ws.Cells.Item(41, 1).Value = "Levy"
$ws.Cells.Item(41, 2).Value = 2014
$ws.Cells.Item(41, 3).Value = 2
$ws.Cells.Item(41, 4).Value = "Medium-textured"
$ws.Cells.Item(41, 5).Value = 4.272553767
$ws.Cells.Item(41, 6).Value = 0.144
$ws.Cells.Item(41, 7).Value = 3
$ws.Cells.Item(41, 8).Value = 2.879058807
$ws.Cells.Item(41, 9).Value = 0.216
$ws.Cells.Item(41, 10).Value = 3

$ws.Cells.Item(42, 1).Value = "Levy"
$ws.Cells.Item(42, 2).Value = 2014
$ws.Cells.Item(42, 3).Value = 2
$ws.Cells.Item(42, 4).Value = "Medium-textured"
$ws.Cells.Item(42, 5).Value = 3.826531262
$ws.Cells.Item(42, 6).Value = 0.144
$ws.Cells.Item(42, 7).Value = 3
$ws.Cells.Item(42, 8).Value = 1.747201758
$ws.Cells.Item(42, 9).Value = 0.144
$ws.Cells.Item(42, 10).Value = 3

$ws.Cells.Item(43, 1).Value = "Levy"
$ws.Cells.Item(43, 2).Value = 2014
$ws.Cells.Item(43, 3).Value = 2
$ws.Cells.Item(43, 4).Value = "Fine-textured"
$ws.Cells.Item(43, 5).Value = 1.816866641
$ws.Cells.Item(43, 6).Value = 0.072
$ws.Cells.Item(43, 7).Value = 3
$ws.Cells.Item(43, 8).Value = 1.747201758
$ws.Cells.Item(43, 9).Value = 0.216
$ws.Cells.Item(43, 10).Value = 3

# Activate the SAR sheet and select the full data range, matching the
# final selection/activeTab state recorded in the workbook.
$ws.Activate()
$ws.Range("A1:J43").Select()
